$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Simple single-value cell updates (row -> old/new value)
$updates = @{
    1  = "0M"
    2  = "0M"
    3  = "0M"
    4  = "1058"
    5  = "0.00002"
    6  = "0.00089"
    7  = "0.00018"
    8  = "0.00006"
    9  = "0.00031"
    10 = "0.00038"
    11 = "0.00051"
    12 = "0.22896"
}

foreach ($row in $updates.Keys) {
    $cell = $t.Cell($row, 1)
    $cell.Range.Text = $updates[$row]
}

# Rows that collapse a multi-run, tab-separated list of values down to a
# single value (stray extra stats removed from README prep step).
$collapse = @{
    44 = "99.83"
    45 = "0.23"
    46 = "136"
}

foreach ($row in $collapse.Keys) {
    $cell = $t.Cell($row, 1)
    $cell.Range.Text = $collapse[$row]
}
